$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The recorded edit re-orders the full content of data rows 2-10: each row's
# whole set of cell values (the same observation record) ends up on a
# different row afterwards. The mapping below gives, for every destination
# row, which row currently holds the data that must be copied there.
$mapping = @{
    2  = 4
    3  = 6
    4  = 2
    5  = 3
    6  = 10
    8  = 5
    10 = 8
}

# A handful of columns are deliberately left untouched for every row
# (regardless of whether that row is a copy destination):
#   - Y and AA always contain the text "2023-09-14" in every data row;
#     writing that text back through a Range value assignment would make
#     Excel auto-coerce it into a serial date, silently changing the cell's
#     stored type/format even though the visible content never changes.
#   - I, K, AT and AY are present in the sheet as explicit-but-empty cells
#     in every data row (no actual value, before or after the edit).
#     Writing an empty value back through COM actually deletes the cell
#     outright instead of keeping it as an empty cell, which would make
#     rows handled by this script diverge from rows left alone. Skipping
#     them keeps every row's treatment of these columns identical.
# The remaining columns are copied in contiguous blocks - A:H, J, L:X, Z,
# AB:AS and AU:AX - i.e. every column except I, K, Y, AA, AT, AY.
$segments = @("A:H", "J:J", "L:X", "Z:Z", "AB:AS", "AU:AX")

# Snapshot every row that is used as a source first, so that writing the
# destination rows never clobbers data that still needs to be read later.
$rowsToRead = $mapping.Values | Sort-Object -Unique
$snapshots = @{}
foreach ($r in $rowsToRead) {
    $bySegment = @{}
    foreach ($seg in $segments) {
        $parts = $seg.Split(":")
        $rangeAddr = "$($parts[0])$r`:$($parts[1])$r"
        $bySegment[$seg] = $ws.Range($rangeAddr).Value()
    }
    $snapshots[$r] = $bySegment
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($seg in $segments) {
        $parts = $seg.Split(":")
        $rangeAddr = "$($parts[0])$destRow`:$($parts[1])$destRow"
        $ws.Range($rangeAddr).Value = $snapshots[$srcRow][$seg]
    }
}
